$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.097.78"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.75%  "
$ws.Range("D3").Value = "'1.906.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.44%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'0.7444"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.65%  "
$ws.Range("D6").Value = "'243.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'0.3088"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.91%  "
$ws.Range("E9").Value = "  -5.43%  "
$ws.Range("D10").Value = "'0.06985"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.62%  "
$ws.Range("D11").Value = "'0.08079"
$ws.Range("D11").Style = "Normal"
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").Value = "'0.7679"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.75%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.958.97"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.17%  "
$ws.Range("D14").Value = "'5.317"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.50%  "
$ws.Range("E15").Value = "  -0.83%  "
$ws.Range("D16").Value = "'14.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.63%  "
$ws.Range("D17").Value = "'30.106.79"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.75%  "
$ws.Range("D18").Value = "'6.084"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("D19").Value = "'0.000007828"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.98%  "
$ws.Range("D20").Value = "'240.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.89%  "
$ws.Range("D21").Value = "'2.161.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.18%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D24").Value = "'7.107"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.99%  "
$ws.Range("D25").Value = "'9.386"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.86%  "
$ws.Range("D26").Value = "'167.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.25%  "
$ws.Range("D27").Value = "'19.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.40%  "
$ws.Range("D28").Value = "'0.1274"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.34%  "
$ws.Range("D29").Value = "'2.050"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.97%  "
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("D31").Value = "'1.351"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.23%  "
$ws.Range("D32").Value = "'4.336"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.31%  "
$ws.Range("E33").Value = "  -1.53%  "
$ws.Range("D34").Value = "'0.05241"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.83%  "
$ws.Range("D35").Value = "'1.306"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.19%  "
$ws.Range("D36").Value = "'0.7477"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.27%  "
$ws.Range("D37").Value = "'2.724"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.15%  "
$ws.Range("D38").Value = "'0.01972"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.55%  "
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("D40").Value = "'6.336"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.90%  "
$ws.Range("D41").Value = "'0.4499"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("D42").Value = "'74.47"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.59%  "
$ws.Range("D43").Value = "'1.975"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.36%  "
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("D45").Value = "'0.8404"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").Value = "'7.743"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.36%  "
$ws.Range("D47").Value = "'101.85"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.17%  "
$ws.Range("D48").Value = "'9.937"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.93%  "
$ws.Range("D49").Value = "'2.086.51"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.48%  "
$ws.Range("D50").Value = "'36.74"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.36%  "
$ws.Range("D51").Value = "'0.1182"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.56%  "
